# Updated with new release 1.6.3 release version
#
# The underlying test data changed the "MessageType" value used for the
# KAFKA verification row (row 3) from "ProtoBuffMessageType" to
# "ProtobufType", and the active selection on the sheet moved from L4 to K15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductCreated-Event")

# K3 holds the MessageType for the "VERIFY_PRODUCT_CREATED_EVENT" (KAFKA) row.
$ws.Range("K3").Value = "ProtobufType"

# Move the active selection to K15 (matches the saved cursor position).
$ws.Range("K15").Select()
